# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" sheet before the current "2022-Q2" sheet by
#    duplicating "2022-Q2" (this carries over all formatting/styles),
#    renaming it, trimming it down to the new data's extent, and
#    overwriting its contents with the 2022-Q3 fund-holding figures.
# 2) Insert a new row into the "总计" (totals) summary sheet with the
#    2022-Q3 aggregate figures, shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Step 1: duplicate "2022-Q2" to create the new "2022-Q3" sheet
# ---------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$sheetQ2.Copy($sheetQ2)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The duplicated sheet has 46 data rows (rows 2-46); the new quarter only
# has 7 fund rows (rows 2-8), so drop the extra rows.
$newSheet.Rows("9:46").Delete()

$fundData = @(
    @(0, "004263", "华安沪港深机会混合",     "8.27", "93.38", "7.16", "0.5921", 3),
    @(1, "012188", "华安优势龙头混合A",       "6.29", "93.24", "5.71", "0.3592", 5),
    @(2, "014539", "华安优势精选混合A",       "1.12", "93.26", "6.39", "0.0716", 5),
    @(3, "006768", "华安沪港深优选混合",       "0.72", "93.08", "6.42", "0.0462", 5),
    @(4, "012189", "华安优势龙头混合C",       "0.62", "93.24", "5.71", "0.0354", 5),
    @(5, "003413", "华泰柏瑞新经济沪港深混合", "0.42", "86.45", "5.27", "0.0221", 9),
    @(6, "014540", "华安优势精选混合C",       "0.24", "93.26", "6.39", "0.0153", 5)
)

$r = 2
foreach ($row in $fundData) {
    $newSheet.Range("A$r").Value = $row[0]

    # B..G are textual (fund code / name / numeric-looking strings that must
    # stay text, matching the source data's inlineStr typing). Force the
    # Text number format first so Excel doesn't re-interpret the numeric-
    # looking values as numbers, then clear the format again afterwards so
    # the cell is left on the default style (only the text *type* needs to
    # stick, not a lingering "@" number format).
    $textRange = $newSheet.Range("B$r`:G$r")
    $textRange.NumberFormat = "@"
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").Value = $row[6]
    $textRange.ClearFormats()

    $newSheet.Range("H$r").Value = $row[7]
    $r++
}

# ---------------------------------------------------------------
# Step 2: insert the 2022-Q3 row into the "总计" summary sheet
# ---------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Rows.Item(2).ClearFormats()

# Re-apply the column-A number style from the row below (keeps the same
# style index the other index cells in column A use).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 1.14

# Column A is a 0-based row index; renumber the rows that got shifted down
# (they kept their old index values after the row insert).
for ($i = 3; $i -le 9; $i++) {
    $totalSheet.Range("A$i").Value = $i - 2
}

Write-Output "2022-Q3 data added"
